$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2-97) forward by 11 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 11
}

# Update retrained model values in column B for the affected rows
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(18, 2).Value = 3
$ws.Cells.Item(19, 2).Value = 3
$ws.Cells.Item(20, 2).Value = 5
$ws.Cells.Item(21, 2).Value = 8
$ws.Cells.Item(22, 2).Value = 94
$ws.Cells.Item(23, 2).Value = 106
$ws.Cells.Item(25, 2).Value = 143
$ws.Cells.Item(26, 2).Value = 427
$ws.Cells.Item(27, 2).Value = 454
$ws.Cells.Item(28, 2).Value = 485
$ws.Cells.Item(29, 2).Value = 520
$ws.Cells.Item(30, 2).Value = 1016
$ws.Cells.Item(31, 2).Value = 1057
$ws.Cells.Item(32, 2).Value = 1100
$ws.Cells.Item(33, 2).Value = 1145
$ws.Cells.Item(34, 2).Value = 1533
$ws.Cells.Item(35, 2).Value = 1570
$ws.Cells.Item(36, 2).Value = 1608
$ws.Cells.Item(37, 2).Value = 1648
$ws.Cells.Item(38, 2).Value = 1940
$ws.Cells.Item(39, 2).Value = 1969
$ws.Cells.Item(40, 2).Value = 1997
$ws.Cells.Item(41, 2).Value = 2019
$ws.Cells.Item(42, 2).Value = 2184
$ws.Cells.Item(43, 2).Value = 2198
$ws.Cells.Item(44, 2).Value = 2213
$ws.Cells.Item(45, 2).Value = 2226
$ws.Cells.Item(46, 2).Value = 2288
$ws.Cells.Item(47, 2).Value = 2296
$ws.Cells.Item(48, 2).Value = 2303
$ws.Cells.Item(49, 2).Value = 2306
$ws.Cells.Item(50, 2).Value = 2309
$ws.Cells.Item(51, 2).Value = 2310
$ws.Cells.Item(52, 2).Value = 2308
$ws.Cells.Item(53, 2).Value = 2303
$ws.Cells.Item(54, 2).Value = 2234
$ws.Cells.Item(55, 2).Value = 2224
$ws.Cells.Item(56, 2).Value = 2211
$ws.Cells.Item(57, 2).Value = 2198
$ws.Cells.Item(58, 2).Value = 2070
$ws.Cells.Item(59, 2).Value = 2049
$ws.Cells.Item(60, 2).Value = 2026
$ws.Cells.Item(61, 2).Value = 2001
$ws.Cells.Item(62, 2).Value = 1789
$ws.Cells.Item(63, 2).Value = 1756
$ws.Cells.Item(64, 2).Value = 1726
$ws.Cells.Item(65, 2).Value = 1695
$ws.Cells.Item(66, 2).Value = 1375
$ws.Cells.Item(67, 2).Value = 1335
$ws.Cells.Item(68, 2).Value = 1300
$ws.Cells.Item(69, 2).Value = 1264
$ws.Cells.Item(70, 2).Value = 835
$ws.Cells.Item(71, 2).Value = 796
$ws.Cells.Item(72, 2).Value = 761
$ws.Cells.Item(73, 2).Value = 726
$ws.Cells.Item(74, 2).Value = 338
$ws.Cells.Item(75, 2).Value = 310
$ws.Cells.Item(76, 2).Value = 285
$ws.Cells.Item(77, 2).Value = 266
$ws.Cells.Item(78, 2).Value = 75
$ws.Cells.Item(79, 2).Value = 61
$ws.Cells.Item(80, 2).Value = 51
$ws.Cells.Item(81, 2).Value = 44
$ws.Cells.Item(82, 2).Value = 5
$ws.Cells.Item(83, 2).Value = 4
$ws.Cells.Item(84, 2).Value = 4
$ws.Cells.Item(85, 2).Value = 4
$ws.Cells.Item(86, 2).Value = 3
$ws.Cells.Item(87, 2).Value = 2
$ws.Cells.Item(88, 2).Value = 2
$ws.Cells.Item(89, 2).Value = 2

Write-Output "Applied Kahraman model retraining updates"
